# Splice in new accelerometer/gyroscope samples:
#   - one new sample inserted right after the header (becomes the new row 2)
#   - nine new samples appended after what is now the last existing row
# The timestamp (A) / label (B) columns are simple, regenerated sequences
# (0, 100, 200, ... ms, all labeled "walkingToRunning") that just grow to
# cover the new row count; only the sensor columns C:H actually shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- existing sensor readings (C:H) for the original 20 data rows (rows 2-21) ---
$existingRange = $ws.Range("C2:H21")
$existing = $existingRange.Value2

$oldRowCount = 20
$cols = 6

# New sample to prepend (becomes row 2)
$newFirst = @(-7.313242021728997, -8.175233658622219, 2.331058922935893, 0.6115283966064453, -0.0169776529073715, -0.1087901294231414)

# New samples to append after the (shifted) original data
$newTail = @(
    ,(7.333738565444815, -26.17261409759575, 7.461308479309277, -11.22835350036621, -15.27582550048828, -1.236372590065002)
    ,(0.3924275192565041, 8.770900960062418, -10.52958610946054, 4.94928503036499, -15.6870174407959, 4.060655117034912)
    ,(-27.41900163538287, -21.30080885045613, 1.806023990406281, 3.654456377029419, -6.942261695861816, 2.724813222885132)
    ,(-26.291865657358, -17.62141748035735, 14.67751483356255, 6.618554592132568, 4.984438896179199, -4.156262397766113)
    ,(-12.81046040852863, 7.074593609454595, 12.13140960768156, -4.189085960388184, 1.281579732894898, 2.08831787109375)
    ,(-5.703531527051739, -11.14156565946689, 10.64924546316558, -3.292665958404541, 1.869073033332825, 3.439073085784912)
    ,(5.73159689061778, -12.49397951013913, 3.91989309647508, -4.698282241821289, 7.67050313949585, -1.46966552734375)
    ,(-10.9989599106356, 0.90863177355609, -9.010537198945535, 1.039232015609741, 13.21467208862305, -9.619471549987791)
    ,(9.205162721520956, -32.48499697329918, -12.15963486129167, 5.204416275024414, -5.015731334686279, -0.8006793856620789)
)

$newRowCount = 1 + $oldRowCount + $newTail.Length

# Build the combined C:H block: new row, then the old rows (shifted down by one), then the appended rows
$block = New-Object 'object[,]' $newRowCount, $cols

for ($c = 0; $c -lt $cols; $c++) { $block[0, $c] = $newFirst[$c] }

for ($r = 1; $r -le $oldRowCount; $r++) {
    for ($c = 0; $c -lt $cols; $c++) {
        $block[$r, $c] = $existing[$r, $c + 1]
    }
}

for ($i = 0; $i -lt $newTail.Length; $i++) {
    $destRow = 1 + $oldRowCount + $i
    for ($c = 0; $c -lt $cols; $c++) {
        $block[$destRow, $c] = $newTail[$i][$c]
    }
}

$lastRow = 1 + $newRowCount  # row 1 is the header, data starts at row 2
$ws.Range("C2:H$lastRow").Value = $block

# --- regenerate the timestamp (A) and label (B) columns for the new row count ---
$ab = New-Object 'object[,]' $newRowCount, 2
for ($r = 0; $r -lt $newRowCount; $r++) {
    $ab[$r, 0] = $r * 100
    $ab[$r, 1] = "walkingToRunning"
}
$ws.Range("A2:B$lastRow").Value = $ab

Write-Output "Spliced in 10 new samples; data now spans rows 2:$lastRow"
